$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Rows("42:42").Insert()
Write-Host "inserted"
